$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-11: new domain (request.biomed.hk) and id incremented by 10.
# Keep column B as text (matches original inline-string text cells, not numbers).
for ($r = 2; $r -le 11; $r++) {
    $oldId = 24082035 + ($r - 2)
    $newId = $oldId + 10
    $ws.Cells.Item($r, 1).Value = "http://request.biomed.hk/12p?name=$newId"
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = "$newId"
}

# Add new rows 12 and 13, continuing the same id sequence.
$ids = @(24082055, 24082056)
$rowIndex = 12
foreach ($id in $ids) {
    $ws.Cells.Item($rowIndex, 1).Value = "http://request.biomed.hk/12p?name=$id"
    $ws.Cells.Item($rowIndex, 2).NumberFormat = "@"
    $ws.Cells.Item($rowIndex, 2).Value = "$id"
    $rowIndex++
}
